$d = $word.ActiveDocument

# Swahili (Kenya) translations for facilitator guideline table/labels.
$replacements = @(
    @("Video Title", "Kichwa cha Video"),
    @("Topic", "Mada"),
    @("Aim(s)", "Malengo"),
    @("Length", "Urefu"),
    @("Camp Location", "Mahali pa Kambi"),
    @("Facilitators", "Wawezeshaji"),
    @("N. of students", "N. ya wanafunzi"),
    @("Date", "Tarehe"),
    @("Resources", "Rasilimali"),
    @("needed", "inahitajika"),
    @("Preparations", "Maandalizi"),
    @("Video time", "Muda wa video"),
    @("What facilitator does", "Mwezeshaji anafanya nini"),
    @("What learners do", "Wanachofanya wanafunzi"),
    @("General VMC Video Introduction", "Utangulizi Mkuu wa Video ya VMC"),
    @("Video Introduction", "Utangulizi wa Video"),
    @("Riddle", "Kitendawili"),
    @("Assist the process, provoke thoughts", "Kusaidia mchakato, kuchochea mawazo"),
    @("Solution", "Suluhisho")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Default document language: Swahili (Tanzania) -> Swahili (Kenya)
$d.Styles(-1).Font.LanguageID = "sw-KE"
